$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for "DOXIRAZOL 60 MG 14 CAPS." after the DIMRA row ---
# Before: row 11 held FUTAPAN; after insert, row 11 is blank and FUTAPAN (and
# everything after it) shifts down to row 12 onward.
$ws.Rows.Item(11).Insert()

# Copy formatting (styles/borders/number-formats) from the row below (the
# shifted-down FUTAPAN row, which still carries the correct data-row style)
# into the newly blank row so the new row matches the rest of the table.
$ws.Range("A12:Q12").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4122)

# Match the row height used by the other data rows in this table.
$ws.Rows.Item(11).RowHeight = 25.5

# Re-create the merged cell regions for the new row (Insert() does not copy
# merges automatically).
$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

# --- Write the new row's data ---
$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "DOXIRAZOL 60 MG 14 CAPS."
$ws.Range("H11").Value = "1:0"
$ws.Range("L11").Value = "1"
$ws.Range("N11").Value = "101.00"
$ws.Range("P11").Value = "101.0000"
$ws.Range("Q11").Value = "1:0"

# --- Update DIMRA row (row 10): new sale price + new order-count ratio ---
$ws.Range("H10").Value = "0:0"
$ws.Range("P10").Value = "105.0000"
$ws.Range("Q10").Value = "1:1"

# --- The remaining rows (12-19) now hold the data that used to sit one row
# higher; rewrite each one explicitly so every cell matches the target file.
$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "FUTAPAN 40 MG VIAL I.V."
$ws.Range("H12").Value = "8:0"
$ws.Range("N12").Value = "59.50"
$ws.Range("P12").Value = "59.5000"
$ws.Range("Q12").Value = "1:0"

$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "MICONAZ 2% ORAL GEL 20 GM"
$ws.Range("H13").Value = "1:0"
$ws.Range("N13").Value = "23.00"
$ws.Range("P13").Value = "23.0000"
$ws.Range("Q13").Value = "1:0"

$ws.Range("A14").Value = 8
$ws.Range("C14").Value = "SILDEN 100 MG 8F.C. TABS"
$ws.Range("H14").Value = "2:7"
$ws.Range("N14").Value = "66.00"
$ws.Range("P14").Value = "66.0000"
$ws.Range("Q14").Value = "1:0"

$ws.Range("A15").Value = 9
$ws.Range("C15").Value = "WATER FOR INJECTION AMP. 5 ML"
$ws.Range("H15").Value = "8447:0"
$ws.Range("L15").Value = "1"
$ws.Range("N15").Value = "2.00"
$ws.Range("P15").Value = "2.0000"
$ws.Range("Q15").Value = "1:0"

$ws.Range("A16").Value = 10
$ws.Range("C16").Value = "سرنجات 3 سم"
$ws.Range("H16").Value = "0:0"
$ws.Range("L16").Value = "0"
$ws.Range("N16").Value = "2.00"
$ws.Range("P16").Value = "2.0000"
$ws.Range("Q16").Value = "1:0"

$ws.Range("A17").Value = 11
$ws.Range("C17").Value = "سرنجات 5 سم"
$ws.Range("H17").Value = "0:0"
$ws.Range("L17").Value = "0"
$ws.Range("N17").Value = "3.00"
$ws.Range("P17").Value = "3.0000"
$ws.Range("Q17").Value = "1:0"

$ws.Range("A18").Value = 12
$ws.Range("C18").Value = "شامبو جونسون 200مللى"
$ws.Range("H18").Value = "1:0"
$ws.Range("L18").Value = "0"
$ws.Range("N18").Value = "50.00"
$ws.Range("P18").Value = "50.0000"
$ws.Range("Q18").Value = "1:0"

# --- New row 19: "كالونا" (previously the last row, now pushed to row 19) ---
$ws.Range("A19").Value = 13
$ws.Range("C19").Value = "كالونا "
$ws.Range("H19").Value = "0:0"
$ws.Range("L19").Value = "0"
$ws.Range("N19").Value = "15.00"
$ws.Range("P19").Value = "15.0000"
$ws.Range("Q19").Value = "1:0"

# --- Update the totals row (now row 20) with the new sum ---
$ws.Range("P20").Value = 520.17999999999995

# --- Update the timestamp in the footer (now row 21) ---
$ws.Range("A21").Value = "Monday, 4 August, 2025 10:31 AM"

Write-Host "edit complete"
